# Fruta / hortaliza, semanal
# Insert a new weekly Cereza (cherry) price record as row 65, shifting all
# subsequent rows (old 65-80) down by one (new 66-81). This matches the
# weekly data-refresh pattern used in this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 65; existing rows 65:80 move to 66:81.
$ws.Rows("65:65").Insert()

# Populate the newly inserted row 65 with this week's record.
$ws.Range("A65").Value = 7
$ws.Range("B65").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C65").Value = "Ñuble"
$ws.Range("D65").Value = 44559
$ws.Range("E65").Value = 16
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100103
$ws.Range("H65").Value = "Frutos de hueso (carozo)"
$ws.Range("I65").Value = 100103001
$ws.Range("J65").Value = "Cereza"
$ws.Range("K65").Value = "Lapins"
$ws.Range("L65").Value = "Primera"
$ws.Range("M65").Value = 120
$ws.Range("N65").Value = 5500
$ws.Range("O65").Value = 6000
$ws.Range("P65").Value = 5750
$ws.Range("Q65").Value = "`$/bandeja 10 kilos"
$ws.Range("R65").Value = "Provincia de Curicó"
$ws.Range("S65").Value = 575
$ws.Range("T65").Value = 10
